$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "NA" values for the duplicate_image_filename column (E)
# for the practice (rows 2-5) and main trial (rows 6-21) rows.
$ws.Range("E2:E21").Value = "NA"
